# "Refined metadata to be additional tab"
#
# 1. Update the time_taken (F column) timestamps on the existing "data" sheet.
# 2. Add a new "metadata" worksheet (after "data") describing the panel pull,
#    re-using the bold/border/centered header style already used on "data".

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- 1. refresh the recorded time_taken values on the "data" sheet ---------
$ws1.Range("F2").Value = "2021-10-05 14:21:22.155191"
$ws1.Range("F3").Value = "2021-10-05 14:21:22.155201"
$ws1.Range("F4").Value = "2021-10-05 14:21:22.155204"
$ws1.Range("F5").Value = "2021-10-05 14:21:22.155207"
$ws1.Range("F6").Value = "2021-10-05 14:21:22.155210"

# --- 2. add the "metadata" worksheet right after "data" --------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "metadata"

# Header row
$ws2.Range("B1").Value = "data_name"
$ws2.Range("C1").Value = "data_id"
$ws2.Range("D1").Value = "data_version"
$ws2.Range("E1").Value = "data_version_created"
$ws2.Range("F1").Value = "panel_query_time"
$ws2.Range("G1").Value = "panel_get_request"

# Data row
$ws2.Range("A2").Value = 0
$ws2.Range("B2").Value = "Leber hereditary optic neuropathy"
$ws2.Range("C2").Value = 530
$ws2.Range("D2").Value = "'1.9"
$ws2.Range("E2").Value = "2021-03-18T16:43:34.412347Z"
$ws2.Range("F2").Value = "2021-10-05 14:21:22.151589"
$ws2.Range("G2").Value = "https://panelapp.genomicsengland.co.uk/api/v1/panels/530/?format=json"

# Re-use the existing header style (bold font, thin border, centered/top
# aligned) from the "data" sheet's header row / index column instead of
# building a brand new style.
$ws1.Range("B1:F1").Copy()
$ws2.Range("B1:F1").PasteSpecial(-4122)
$ws1.Range("F1").Copy()
$ws2.Range("G1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws2.Range("A2").PasteSpecial(-4122)

# D2 ("1.9") must stay plain (unstyled) text like the source data -- reset
# its style back to the default after the quote-prefix forced text type.
$ws1.Range("D2").Copy()
$ws2.Range("D2").PasteSpecial(-4122)

Write-Output "metadata sheet added"
